$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the existing header style (s=1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data columns I (I0) and J (IF) for rows 2-74
$data = @(
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(5, 5),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(5, 6),
    @(11, 11),
    @(6, 6),
    @(5, 5),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(2, 3),
    @(6, 6),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(10, 10),
    @(7, 7),
    @(10, 10),
    @(6, 6),
    @(6, 7),
    @(7, 7),
    @(8, 8),
    @(10, 10),
    @(5, 6),
    @(3, 4),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(4, 5),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(5, 5),
    @(5, 5),
    @(6, 6),
    @(3, 3),
    @(6, 6),
    @(6, 6)
)

for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $row = 2 + $idx
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}
